$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New "Sheet3" - Allergies list
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "Sheet3"

$ws3.Range("A1").Value = "Allergies"
$ws3.Range("A1").Font.Bold = $true
$ws3.Range("A2").Value = "milk"
$ws3.Range("A3").Value = "soy"
$ws3.Range("A4").Value = "egg"
$ws3.Range("A5").Value = "sesame"
$ws3.Range("A6").Value = "shellfish"
$ws3.Range("A7").Value = "seafood"

# ---------------------------------------------------------------------------
# New "Sheet4" - Nut Allergies list
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "Sheet4"

$ws4.Range("A1").Value = "Nut Allergies"
$ws4.Range("A1").Font.Bold = $true
$ws4.Range("A2").Value = "peanuts"
$ws4.Range("A3").Value = "walnuts"
$ws4.Range("A4").Value = "almond"
$ws4.Range("A5").Value = "hazelnut"
$ws4.Range("A6").Value = "cashew"
$ws4.Range("A7").Value = "pecan"
$ws4.Range("A8").Value = "pistachio"

# ---------------------------------------------------------------------------
# Selections matching final saved state, and make Sheet4 the active tab
# ---------------------------------------------------------------------------
[void]$ws3.Range("A1:A7").Select()
[void]$ws4.Range("A1:A8").Select()
$ws4.Activate()
